$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Overview sheet: "In Translation" -> "Handed back: in sync with en-US"
#    (E2, F2, E3, F3 all share this string; the Status cells on the
#    zh-cn/de-de sheets -- C2/C3 -- reference the same shared string,
#    so updating the text anywhere it is used keeps them all in sync.)
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the zh-cn / de-de columns on the Overview sheet to fit the
# longer status text.
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777050018311
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777050018311

# ------------------------------------------------------------------
# 2) zh-cn sheet: fill in "Latest Target File", "Latest Handback File"
#    and "Latest Handback DateTime" for the two rows (report generated
#    for the handback).
# ------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

$wsZhCn.Range("I2").Value = "31a62b58-ca2f-4e96-b1f2-b033c8f9014e.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48785027a6e8d50bab2e3c297e1777d7192b3333/e2e/31a62b58-ca2f-4e96-b1f2-b033c8f9014e.md", "", "", "31a62b58-ca2f-4e96-b1f2-b033c8f9014e.md") | Out-Null
$wsZhCn.Range("J2").Value = "31a62b58-ca2f-4e96-b1f2-b033c8f9014e.7825fabcd8d28dad0e6c213f0f00fca314026d00.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-06 11:35:32"

$wsZhCn.Range("I3").Value = "687502d2-81fd-4b8c-b7ab-b5a49e3d6ad0.md"
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48785027a6e8d50bab2e3c297e1777d7192b3333/e2e/687502d2-81fd-4b8c-b7ab-b5a49e3d6ad0.md", "", "", "687502d2-81fd-4b8c-b7ab-b5a49e3d6ad0.md") | Out-Null
$wsZhCn.Range("J3").Value = "687502d2-81fd-4b8c-b7ab-b5a49e3d6ad0.99e4d2f779a705983d1d740a68af9157ae828396.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-06 11:35:32"

$wsZhCn.Columns.Item(3).ColumnWidth = 29.9777050018311
$wsZhCn.Columns.Item(9).ColumnWidth = 40
$wsZhCn.Columns.Item(10).ColumnWidth = 40

# ------------------------------------------------------------------
# 3) de-de sheet: same updates, different handback timestamp
# ------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDeDe.Range("I2").Value = "31a62b58-ca2f-4e96-b1f2-b033c8f9014e.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48785027a6e8d50bab2e3c297e1777d7192b3333/e2e/31a62b58-ca2f-4e96-b1f2-b033c8f9014e.md", "", "", "31a62b58-ca2f-4e96-b1f2-b033c8f9014e.md") | Out-Null
$wsDeDe.Range("J2").Value = "31a62b58-ca2f-4e96-b1f2-b033c8f9014e.7825fabcd8d28dad0e6c213f0f00fca314026d00.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-06 11:35:50"

$wsDeDe.Range("I3").Value = "687502d2-81fd-4b8c-b7ab-b5a49e3d6ad0.md"
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/48785027a6e8d50bab2e3c297e1777d7192b3333/e2e/687502d2-81fd-4b8c-b7ab-b5a49e3d6ad0.md", "", "", "687502d2-81fd-4b8c-b7ab-b5a49e3d6ad0.md") | Out-Null
$wsDeDe.Range("J3").Value = "687502d2-81fd-4b8c-b7ab-b5a49e3d6ad0.99e4d2f779a705983d1d740a68af9157ae828396.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-06 11:35:50"

$wsDeDe.Columns.Item(3).ColumnWidth = 29.9777050018311
$wsDeDe.Columns.Item(9).ColumnWidth = 40
$wsDeDe.Columns.Item(10).ColumnWidth = 40
